$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New transaction rows (22/23 Jan 2026 bank statement entries) ---------
$newRows = @(
    @("23/01/2026", "23/01/2026", "673 GADIS LA PALLOZA     A CORUÑA     ES", "PAGO CON TARJETA", -62.55, "EUR", 1155.39, "Hogar", "Supermercado"),
    @("18/01/2026", "23/01/2026", "Macenda (-A Coruña  17:38", "PAGO CON TARJETA", -2.45, "EUR", 1217.94, "Coche", "Peaje"),
    @("18/01/2026", "23/01/2026", "A Coruña -R.Guísamo 12:25", "PAGO CON TARJETA", -2.05, "EUR", 1220.39, "Coche", "Peaje"),
    @("22/01/2026", "22/01/2026", "EL CORTE INGLES-DEPARTAM.", "PAGO CON TARJETA", -44.94, "EUR", 1222.44, "Niñas", "Ropa"),
    @("20/01/2026", "22/01/2026", "EasyPark Espana S.L.U.", "PAGO CON TARJETA", -0.9, "EUR", 1267.38, "Coche", "Parking"),
    @("21/01/2026", "22/01/2026", "EasyPark Espana S.L.U.", "PAGO CON TARJETA", -0.98, "EUR", 1268.28, "Coche", "Parking"),
    @("21/01/2026", "21/01/2026", "673 GADIS LA PALLOZA", "PAGO CON TARJETA", -28.35, "EUR", 1269.26, "Hogar", "Supermercado"),
    @("20/01/2026", "20/01/2026", "673 GADIS LA PALLOZA", "PAGO CON TARJETA", -2.18, "EUR", 1297.61, "Hogar", "Supermercado"),
    @("16/01/2026", "19/01/2026", "APPLE.COM/BILL", "PAGO CON TARJETA", -10.99, "EUR", 1299.79, "Hogar", "Ocio"),
    @("18/01/2026", "19/01/2026", "Netflix.com", "PAGO CON TARJETA", -13.99, "EUR", 1310.78, "Hogar", "Ocio"),
    @("15/01/2026", "19/01/2026", "EasyPark Espana S.L.U.", "PAGO CON TARJETA", -0.69, "EUR", 1324.77, "Coche", "Parking"),
    @("16/01/2026", "19/01/2026", "EasyPark Espana S.L.U.", "PAGO CON TARJETA", -0.98, "EUR", 1325.46, "Coche", "Parking"),
    @("18/01/2026", "19/01/2026", "EasyPark Espana S.L.U.", "PAGO CON TARJETA", -2.59, "EUR", 1326.44, "Coche", "Parking"),
    @("17/01/2026", "19/01/2026", "EasyPark Espana S.L.U.", "PAGO CON TARJETA", -2.54, "EUR", 1329.03, "Coche", "Parking"),
    @("17/01/2026", "19/01/2026", "673 GADIS LA PALLOZA", "PAGO CON TARJETA", -21.74, "EUR", 1331.57, "Hogar", "Supermercado"),
    @("17/01/2026", "19/01/2026", "673 GADIS LA PALLOZA", "PAGO CON TARJETA", -10.73, "EUR", 1353.31, "Hogar", "Supermercado"),
)


$startRow = 37
$endRow = $startRow + $newRows.Count - 1

# Bring the new rows' formatting in line with the existing data rows: copy
# the format (borders/font) of the last existing data row (36) down across
# the whole new block, then fill in the values cell by cell.
$ws.Range("A36:I36").Copy()
$ws.Range(("A" + $startRow + ":I" + $endRow)).PasteSpecial(-4122)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}

# --- Update the view state to match what was selected/scrolled to after
# --- entering the new rows (matches the author's final cursor position).
$win = $excel.ActiveWindow
$win.ScrollRow = 34
$win.ScrollColumn = 3
$ws.Range("I39").Select()

Write-Output ("Added rows " + $startRow + " to " + $endRow)
